# Cập nhật data.xlsx từ công cụ QR
# Inserts a newly-scanned location record as row 2 of the "locations" sheet,
# pushing the existing rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("locations")

# Shift existing data rows (old row 2..5) down by one, keeping them intact.
$ws.Rows.Item(2).Insert()

# Fill in the new row with the freshly scanned record.
$ws.Cells.Item(2, 1).Value  = "wajj1kvqv4r"
$ws.Cells.Item(2, 2).Value  = "y302zsdg"
$ws.Cells.Item(2, 3).Value  = "C"
$ws.Cells.Item(2, 4).Value  = "Phường Tăng Nhơn Phú, Ho Chi Minh City, 71300, Vietnam"
$ws.Cells.Item(2, 5).Value  = "https://www.google.com/maps/search/?api=1&query=10.839061,106.792777"
$ws.Cells.Item(2, 6).Value  = "2025-08-22T09:48:08.014Z"
$ws.Cells.Item(2, 7).Value  = ""
$ws.Cells.Item(2, 8).Value  = ""
$ws.Cells.Item(2, 9).Value  = ""
$ws.Cells.Item(2, 10).Value = ""
$ws.Cells.Item(2, 11).Value = ""
$ws.Cells.Item(2, 12).Value = "C"
$ws.Cells.Item(2, 13).Value = "5c8078db067e40cd"
$ws.Cells.Item(2, 14).Value = "85bc7e5f757639fb1c6a791abcf143e1a4bfadd181e61e5df85b926bd4928ffc"
